$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1162.3
$ws.Range("I19").Value = 1125.8
$ws.Range("K19").Value = 1125.8
$ws.Range("M19").Value = -950.8

$ws.Range("H28").Value = 385.51852
$ws.Range("I28").Value = 343.72726
$ws.Range("K28").Value = 343.72726
$ws.Range("M28").Value = 141.27274

$ws.Range("H43").Value = 3498
$ws.Range("I43").Value = 1001
$ws.Range("K43").Value = 1001
$ws.Range("M43").Value = -932

$ws.Range("H58").Value = 2170.2856
$ws.Range("I58").Value = 1025
$ws.Range("J58").Value = 3029.25
$ws.Range("K58").Value = 3075
$ws.Range("L58").Value = 9087.75
$ws.Range("M58").Value = -2925
$ws.Range("N58").Value = -9387.75

$ws.Range("H98").Value = 764.6
$ws.Range("J98").Value = 618
$ws.Range("L98").Value = 618
$ws.Range("N98").Value = -3614

$ws.Range("H122").Value = 764.6
$ws.Range("J122").Value = 618
$ws.Range("L122").Value = 1854
$ws.Range("N122").Value = -6754

$ws.Range("H129").Value = 1533.1666
$ws.Range("I129").Value = 1300
$ws.Range("J129").Value = 1999.5
$ws.Range("K129").Value = 3900
$ws.Range("L129").Value = 5998.5
$ws.Range("M129").Value = 1100
$ws.Range("N129").Value = -15998.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2859
$ws.Range("I45").Value = 1955.1666
$ws.Range("K45").Value = 1955.1666
$ws.Range("M45").Value = -1578.1666

$ws.Range("H110").Value = 1857.579
$ws.Range("I110").Value = 1661.875
$ws.Range("J110").Value = 1999.909
$ws.Range("K110").Value = 1661.875
$ws.Range("L110").Value = 1999.909
$ws.Range("M110").Value = 383.125
$ws.Range("N110").Value = -6089.909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37037.445
$ws.Range("J81").Value = 40417.125
$ws.Range("L81").Value = 40417.125
$ws.Range("N81").Value = -42539.125

$ws.Range("H84").Value = 37037.445
$ws.Range("J84").Value = 40417.125
$ws.Range("L84").Value = 121251.375
$ws.Range("N84").Value = -131859.375

$ws.Range("H99").Value = 13937.375
$ws.Range("I99").Value = 13937.375
$ws.Range("K99").Value = 13937.375
$ws.Range("M99").Value = -12439.375

$ws.Range("H107").Value = 1709.5
$ws.Range("I107").Value = 1455.3846
$ws.Range("K107").Value = 1455.3846
$ws.Range("M107").Value = 464.6153999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2723.0908
$ws.Range("I31").Value = 2723.0908
$ws.Range("K31").Value = 2723.0908
$ws.Range("M31").Value = -2428.0908

$ws.Range("H34").Value = 2723.0908
$ws.Range("I34").Value = 2723.0908
$ws.Range("K34").Value = 2723.0908
$ws.Range("M34").Value = -2521.0908

$ws.Range("H86").Value = 85648.75
$ws.Range("J86").Value = 18316.166
$ws.Range("L86").Value = 18316.166
$ws.Range("N86").Value = -20562.166

$ws.Range("H89").Value = 85648.75
$ws.Range("J89").Value = 18316.166
$ws.Range("L89").Value = 91580.83
$ws.Range("N89").Value = -102812.83

$ws.Range("H105").Value = 16499.889
$ws.Range("I105").Value = 18187.375
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 18187.375
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -16440.375
$ws.Range("N105").Value = -6494

$ws.Range("H122").Value = 9278.058
$ws.Range("I122").Value = 2093.0715
$ws.Range("K122").Value = 6279.2145
$ws.Range("M122").Value = -3829.2145

$ws.Range("H134").Value = 3727.077
$ws.Range("I134").Value = 3643.1667
$ws.Range("J134").Value = 3799
$ws.Range("K134").Value = 10929.5001
$ws.Range("L134").Value = 11397
$ws.Range("M134").Value = -8394.500100000001
$ws.Range("N134").Value = -16467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 12065.5
$ws.Range("I94").Value = 6196.5
$ws.Range("K94").Value = 18589.5
$ws.Range("M94").Value = -17913.5

$ws.Range("H113").Value = 25464.5
$ws.Range("I113").Value = 372.75
$ws.Range("J113").Value = 38010.375
$ws.Range("K113").Value = 1118.25
$ws.Range("L113").Value = 114031.125
$ws.Range("M113").Value = 1051.75
$ws.Range("N113").Value = -118371.125

$ws.Range("H116").Value = 5781.125
$ws.Range("I116").Value = 1562.25
$ws.Range("K116").Value = 4686.75
$ws.Range("M116").Value = -1244.75

$ws.Range("H122").Value = 5379677.5
$ws.Range("I122").Value = 16129032
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 145161288
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -145158838
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17529.21
$ws.Range("I70").Value = 21179.615
$ws.Range("J70").Value = 9620
$ws.Range("K70").Value = 21179.615
$ws.Range("L70").Value = 9620
$ws.Range("M70").Value = -20909.615
$ws.Range("N70").Value = -10160

$ws.Range("H73").Value = 17529.21
$ws.Range("I73").Value = 21179.615
$ws.Range("J73").Value = 9620
$ws.Range("K73").Value = 21179.615
$ws.Range("L73").Value = 9620
$ws.Range("M73").Value = -20243.615
$ws.Range("N73").Value = -11492

$ws.Range("H102").Value = 2752.6667
$ws.Range("J102").Value = 1499.5
$ws.Range("L102").Value = 1499.5
$ws.Range("N102").Value = -4743.5

$ws.Range("H113").Value = 2005.6786
$ws.Range("I113").Value = 1895.9546
$ws.Range("J113").Value = 2408
$ws.Range("K113").Value = 1895.9546
$ws.Range("L113").Value = 2408
$ws.Range("M113").Value = 274.0454
$ws.Range("N113").Value = -6748

$ws.Range("H122").Value = 2486.52
$ws.Range("I122").Value = 1613.5834
$ws.Range("J122").Value = 3292.3076
$ws.Range("K122").Value = 4840.7502
$ws.Range("L122").Value = 9876.9228
$ws.Range("M122").Value = -2390.7502
$ws.Range("N122").Value = -14776.9228

$ws.Range("H126").Value = 2810.4546
$ws.Range("I126").Value = 2600.7144
$ws.Range("K126").Value = 7802.1432
$ws.Range("M126").Value = -5332.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3019
$ws.Range("J7").Value = 4997.5
$ws.Range("L7").Value = 4997.5
$ws.Range("N7").Value = -5221.5

$ws.Range("H40").Value = 26875
$ws.Range("I40").Value = 34166.668
$ws.Range("K40").Value = 34166.668
$ws.Range("M40").Value = -34030.668

$ws.Range("H82").Value = 1877.9333
$ws.Range("I82").Value = 1791.9
$ws.Range("J82").Value = 2050
$ws.Range("K82").Value = 1791.9
$ws.Range("L82").Value = 2050
$ws.Range("M82").Value = -1430.9
$ws.Range("N82").Value = -2772

$ws.Range("H85").Value = 1877.9333
$ws.Range("I85").Value = 1791.9
$ws.Range("J85").Value = 2050
$ws.Range("K85").Value = 1791.9
$ws.Range("L85").Value = 2050
$ws.Range("M85").Value = -543.9000000000001
$ws.Range("N85").Value = -4546

$ws.Range("H126").Value = 3019
$ws.Range("J126").Value = 4997.5
$ws.Range("L126").Value = 14992.5
$ws.Range("N126").Value = -19932.5

$ws.Range("H132").Value = 3025.08
$ws.Range("I132").Value = 2035.4615
$ws.Range("J132").Value = 4097.1665
$ws.Range("K132").Value = 6106.3845
$ws.Range("L132").Value = 12291.4995
$ws.Range("M132").Value = -3576.3845
$ws.Range("N132").Value = -17351.4995

$ws.Range("H136").Value = 11281.083
$ws.Range("I136").Value = 5528.143
$ws.Range("K136").Value = 16584.429
$ws.Range("M136").Value = -14034.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8154.3687
$ws.Range("I81").Value = 2826.4614
$ws.Range("J81").Value = 19698.166
$ws.Range("K81").Value = 5652.9228
$ws.Range("L81").Value = 39396.332
$ws.Range("M81").Value = -4591.9228
$ws.Range("N81").Value = -41518.332

$ws.Range("H84").Value = 8154.3687
$ws.Range("I84").Value = 2826.4614
$ws.Range("J84").Value = 19698.166
$ws.Range("K84").Value = 28264.614
$ws.Range("L84").Value = 196981.66
$ws.Range("M84").Value = -22960.614
$ws.Range("N84").Value = -207589.66

$ws.Range("H107").Value = 1455.1714
$ws.Range("I107").Value = 1119.5
$ws.Range("J107").Value = 1902.7333
$ws.Range("K107").Value = 3358.5
$ws.Range("L107").Value = 5708.199900000001
$ws.Range("M107").Value = -1438.5
$ws.Range("N107").Value = -9548.1999

$ws.Range("H113").Value = 6433
$ws.Range("I113").Value = 4899.5
$ws.Range("J113").Value = 9500
$ws.Range("K113").Value = 14698.5
$ws.Range("L113").Value = 28500
$ws.Range("M113").Value = -12528.5
$ws.Range("N113").Value = -32840

$ws.Range("H122").Value = 84784.29
$ws.Range("I122").Value = 3995.5557
$ws.Range("J122").Value = 230204
$ws.Range("K122").Value = 11986.6671
$ws.Range("L122").Value = 690612
$ws.Range("M122").Value = -9536.667099999999
$ws.Range("N122").Value = -695512

$ws.Range("H136").Value = 7429.154
$ws.Range("I136").Value = 7964.9165
$ws.Range("K136").Value = 23894.7495
$ws.Range("M136").Value = -21344.7495
